$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.654.27'
$ws.Range("E2").Value = '  +1.31%  '
$ws.Range("D3").Value = '1.633.69'
$ws.Range("E3").Value = '  +1.44%  '
$ws.Range("E4").Value = '  -0.09%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '213.09'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.495'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +1.56%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  +1.19%  '
$ws.Range("E9").Value = '  +1.65%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '19.03'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +2.62%  '
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("D12").Value = '1.861.30'
$ws.Range("E12").Value = '  +1.42%  '
$ws.Range("D13").Value = '1.627.69'
$ws.Range("E13").Value = '  +0.14%  '
$ws.Range("E14").Value = '  +1.67%  '
$ws.Range("E15").Value = '  +2.24%  '
$ws.Range("D16").Value = '26.654.42'
$ws.Range("E16").Value = '  +1.40%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '63.21'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +1.49%  '
$ws.Range("D18").Value = '0.0₃0740'
$ws.Range("E18").Value = '  +1.78%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '210.57'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +4.33%  '
$ws.Range("E20").Value = '  -0.09%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '4.31'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +0.89%  '
$ws.Range("E22").Value = '  +1.24%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '6.22'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +3.06%  '
$ws.Range("E24").Value = '  +1.99%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '147.24'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +2.45%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  -0.65%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '6.89'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +4.81%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '15.40'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.98%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.0523'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +5.22%  '
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("E32").Value = '  +1.68%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '2.95'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("E34").Value = '  +1.13%  '
$ws.Range("E35").Value = '  -0.84%  '
$ws.Range("E36").Value = '  +2.43%  '
$ws.Range("D37").Value = '1.170.78'
$ws.Range("E37").Value = '  +0.69%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.812'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  -0.07%  '
$ws.Range("E40").Value = '  +1.90%  '
$ws.Range("E41").Value = '  -0.30%  '
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").Value = '1.770.42'
$ws.Range("E44").Value = '  +1.40%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '92.52'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("E46").Value = '  +1.11%  '
$ws.Range("E47").Value = '  +1.49%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0512'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.82%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '7.56'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  +4.32%  '
$ws.Range("E50").Value = '  +0.26%  '
$ws.Range("E51").Value = '  -0.12%  '
